# Updates cryptos list values (price/volume columns) per the GitHub Actions
# scrape refresh, and swaps the Kaspa/LidoDAOToken row contents (rows 36-37).
# Numeric-looking text values are prefixed with a leading apostrophe so Excel
# keeps them stored as text (matching the source data, e.g. "6.80", "0.102")
# instead of silently coercing them into numbers and dropping trailing zeros.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.911.85"
$ws.Range("E2").Value = "  -1.14%  "
$ws.Range("D3").Value = "2.192.44"
$ws.Range("E3").Value = "  -2.45%  "
$ws.Range("E4").Value = "  +0.08%  "
$ws.Range("D5").Value = "'294.27"
$ws.Range("E5").Value = "  -4.46%  "
$ws.Range("D6").Value = "'88.79"
$ws.Range("E6").Value = "  -6.38%  "
$ws.Range("D7").Value = "'0.565"
$ws.Range("E7").Value = "  -1.18%  "
$ws.Range("D9").Value = "'0.482"
$ws.Range("E9").Value = "  -8.55%  "
$ws.Range("D10").Value = "'31.96"
$ws.Range("E10").Value = "  -8.74%  "
$ws.Range("D11").Value = "'0.0770"
$ws.Range("E11").Value = "  -5.14%  "
$ws.Range("E12").Value = "  -1.34%  "
$ws.Range("D13").Value = "'6.80"
$ws.Range("E13").Value = "  -5.72%  "
$ws.Range("D14").Value = "2.529.14"
$ws.Range("E14").Value = "  -2.40%  "
$ws.Range("D15").Value = "2.258.63"
$ws.Range("E15").Value = "  -4.26%  "
$ws.Range("D16").Value = "'13.10"
$ws.Range("E16").Value = "  -4.47%  "
$ws.Range("D17").Value = "'0.771"
$ws.Range("E17").Value = "  -8.35%  "
$ws.Range("D18").Value = "43.605.46"
$ws.Range("E18").Value = "  -1.21%  "
$ws.Range("D19").Value = "0.0₃0886"
$ws.Range("E19").Value = "  -8.38%  "
$ws.Range("D20").Value = "'5.82"
$ws.Range("E20").Value = "  -9.00%  "
$ws.Range("D21").Value = "'10.75"
$ws.Range("E21").Value = "  -12.63%  "
$ws.Range("D22").Value = "'63.03"
$ws.Range("E22").Value = "  -4.31%  "
$ws.Range("D23").Value = "'232.31"
$ws.Range("D24").Value = "'2.75"
$ws.Range("E24").Value = "  -9.31%  "
$ws.Range("E25").Value = "  +0.73%  "
$ws.Range("D26").Value = "'1.83"
$ws.Range("E26").Value = "  -9.47%  "
$ws.Range("D27").Value = "'2.23"
$ws.Range("E27").Value = "  +0.61%  "
$ws.Range("D28").Value = "'36.15"
$ws.Range("E28").Value = "  -6.11%  "
$ws.Range("D29").Value = "'9.22"
$ws.Range("E29").Value = "  -6.65%  "
$ws.Range("D30").Value = "'19.19"
$ws.Range("E30").Value = "  -4.51%  "
$ws.Range("D31").Value = "'148.15"
$ws.Range("E31").Value = "  -3.91%  "
$ws.Range("D32").Value = "'5.26"
$ws.Range("E32").Value = "  -12.02%  "
$ws.Range("D33").Value = "'2.51"
$ws.Range("E33").Value = "  -5.24%  "
$ws.Range("D34").Value = "'0.0731"
$ws.Range("E34").Value = "  -8.84%  "
$ws.Range("D35").Value = "'0.115"
$ws.Range("E35").Value = "  -3.86%  "
$ws.Range("B36").Value = "Kaspa"
$ws.Range("C36").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D36").Value = "'0.102"
$ws.Range("E36").Value = "  -6.02%  "
$ws.Range("B37").Value = "LidoDAOToken"
$ws.Range("C37").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D37").Value = "'2.80"
$ws.Range("E37").Value = "  -10.14%  "
$ws.Range("D38").Value = "'1.63"
$ws.Range("E38").Value = "  -9.70%  "
$ws.Range("D39").Value = "'0.0282"
$ws.Range("E39").Value = "  -7.54%  "
$ws.Range("D40").Value = "'3.49"
$ws.Range("E40").Value = "  -8.67%  "
$ws.Range("D41").Value = "'3.06"
$ws.Range("E41").Value = "  -12.35%  "
$ws.Range("E42").Value = "  -0.18%  "
$ws.Range("D43").Value = "'12.89"
$ws.Range("E43").Value = "  -12.93%  "
$ws.Range("D44").Value = "1.789.29"
$ws.Range("E44").Value = "  +2.84%  "
$ws.Range("D45").Value = "'1.65"
$ws.Range("E45").Value = "  +3.15%  "
$ws.Range("E46").Value = "  +11.54%  "
$ws.Range("D47").Value = "'0.173"
$ws.Range("E47").Value = "  -10.61%  "
$ws.Range("D48").Value = "'72.29"
$ws.Range("E48").Value = "  -10.36%  "
$ws.Range("D49").Value = "'91.63"
$ws.Range("E49").Value = "  -8.45%  "
$ws.Range("D50").Value = "'64.76"
$ws.Range("E50").Value = "  -8.86%  "
$ws.Range("D51").Value = "2.413.12"
$ws.Range("E51").Value = "  -2.36%  "
